# Reorder rows 12-14 (cyclic rotation) and swap rows 18-19, per the
# upstream re-sort. Only the cells whose value actually changes are
# touched; Y/AA (both "2026-02-16" everywhere here) are left alone so
# Excel does not reinterpret the literal text as a date serial.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 ---
$ws.Range("A12").Value2 = 131187835
$ws.Range("B12").Value2 = 57073
$ws.Range("D12").Value2 = "LC"
$ws.Range("E12").Value2 = 100138
$ws.Range("F12").Value2 = "Tjäder"
$ws.Range("G12").Value2 = "Tetrao urogallus"
$ws.Range("H12").Value2 = "Linnaeus, 1758"
$ws.Range("J12").ClearContents() | Out-Null
$ws.Range("L12").ClearContents() | Out-Null
$ws.Range("M12").Value2 = "äldre spår"
$ws.Range("Q12").Value2 = 511382
$ws.Range("R12").Value2 = 6697458
$ws.Range("S12").Value2 = 25
$ws.Range("Z12").ClearContents() | Out-Null
$ws.Range("AB12").ClearContents() | Out-Null
$ws.Range("AC12").Value2 = "Betad tallkrona."
$ws.Range("AF12").ClearContents() | Out-Null
$ws.Range("AW12").Value2 = "Anna-Lena Thommson"
$ws.Range("AX12").Value2 = "Anna-Lena Thommson"

# --- Row 13 ---
$ws.Range("A13").Value2 = 131191949
$ws.Range("P13").Value2 = "Svartå, Dlr"
$ws.Range("Q13").Value2 = 511393
$ws.Range("R13").Value2 = 6697824
$ws.Range("S13").Value2 = 10
$ws.Range("Z13").Value2 = "10:33"
$ws.Range("AB13").Value2 = "10:33"
$ws.Range("AC13").ClearContents() | Out-Null
$ws.Range("AW13").Value2 = "Lars-Erik Nilsson"
$ws.Range("AX13").Value2 = "Lars-Erik Nilsson, Anna-Lena Thommson"

# --- Row 14 ---
$ws.Range("A14").Value2 = 131187780
$ws.Range("B14").Value2 = 79245
$ws.Range("D14").Value2 = "NT"
$ws.Range("E14").Value2 = 6425
$ws.Range("F14").Value2 = "Garnlav"
$ws.Range("G14").Value2 = "Alectoria sarmentosa"
$ws.Range("H14").Value2 = "(Ach.) Ach."
$ws.Range("J14").ClearContents() | Out-Null
$ws.Range("L14").ClearContents() | Out-Null
$ws.Range("M14").ClearContents() | Out-Null
$ws.Range("P14").Value2 = "Svatå, Dlr"
$ws.Range("Q14").Value2 = 511335
$ws.Range("R14").Value2 = 6697864
$ws.Range("AC14").Value2 = "På gran."
$ws.Range("AF14").ClearContents() | Out-Null

# --- Row 18 ---
$ws.Range("A18").Value2 = 131187762
$ws.Range("B18").Value2 = 79245
$ws.Range("D18").Value2 = "NT"
$ws.Range("E18").Value2 = 6425
$ws.Range("F18").Value2 = "Garnlav"
$ws.Range("G18").Value2 = "Alectoria sarmentosa"
$ws.Range("H18").Value2 = "(Ach.) Ach."
$ws.Range("J18").ClearContents() | Out-Null
$ws.Range("L18").ClearContents() | Out-Null
$ws.Range("M18").ClearContents() | Out-Null
$ws.Range("P18").Value2 = "Svartå, Dlr"
$ws.Range("Q18").Value2 = 511511
$ws.Range("R18").Value2 = 6697866
$ws.Range("AC18").Value2 = "På äldre tall."
$ws.Range("AF18").ClearContents() | Out-Null

# --- Row 19 ---
$ws.Range("A19").Value2 = 131187791
$ws.Range("B19").Value2 = 57073
$ws.Range("D19").Value2 = "LC"
$ws.Range("E19").Value2 = 100138
$ws.Range("F19").Value2 = "Tjäder"
$ws.Range("G19").Value2 = "Tetrao urogallus"
$ws.Range("H19").Value2 = "Linnaeus, 1758"
$ws.Range("J19").ClearContents() | Out-Null
$ws.Range("L19").ClearContents() | Out-Null
$ws.Range("M19").Value2 = "färsk spillning"
$ws.Range("P19").Value2 = "Svatå, Dlr"
$ws.Range("Q19").Value2 = 511301
$ws.Range("R19").Value2 = 6697864
$ws.Range("AC19").ClearContents() | Out-Null
$ws.Range("AF19").ClearContents() | Out-Null
